# TC19_Verify_ShippingMethod.xlsx — "new changes in release part 5"
#
# A new test step was inserted into the keyword-driven test-case sheet,
# right before the existing "CLICK_PRE_ENTERTEXT / ShippingCity" step
# (which, together with everything after it, shifts down by one row).
# The new step is a bare "TINY_SCROLL_DOWN" keyword row (Object /
# ObjectType / Data_descriptor left blank), matching the same kind of
# scroll step that already exists later in the sheet (originally row 48).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
[void]$ws.Activate()

# Push row 42 (and everything below it) down by inserting a fresh row.
[void]$ws.Rows.Item(42).Insert()

# Populate the new row 42 with the new keyword step.
$ws.Cells.Item(42, 1).Value = $null
$ws.Cells.Item(42, 2).Value = "TINY_SCROLL_DOWN"
$ws.Cells.Item(42, 3).Value = $null
$ws.Cells.Item(42, 4).Value = $null
$ws.Cells.Item(42, 5).Value = $null

# Match the bordered-table formatting used by every other data row on
# this sheet (thin border all around each of the 5 columns).
for ($col = 1; $col -le 5; $col++) {
    $cell = $ws.Cells.Item(42, $col)
    $cell.Borders.Item(1).LineStyle = 1
    $cell.Borders.Item(2).LineStyle = 1
    $cell.Borders.Item(3).LineStyle = 1
    $cell.Borders.Item(4).LineStyle = 1
}

# Leave the selection on the newly inserted cell, as in the saved file.
[void]$ws.Range("B42").Select()
